$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.838.23"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.757.33"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "328.26"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  +0.04%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4576"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3497"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "41.97"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07352"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -0.09%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "20.63"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.982"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.170"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.756.98"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "91.73"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.61%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001053"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06408"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -1.38%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.740"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "27.867.86"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.160"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.70%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "162.32"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.04"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.959.78"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.17%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.160"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.61%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "123.19"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.98%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.083"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09297"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.639"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.46%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.539"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "11.74"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06103"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02248"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2065"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.887"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6188"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -1.31%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "7.794"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.06"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  +0.14%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5797"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "122.27"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.925"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.120"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06780"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "72.12"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
